# Add the new leaderboard entry for Carlos Salomao (row 14) and refresh the
# "posicao" ranking + point totals for the rows pushed down beneath it.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Insert a fresh row at position 14 (everything below shifts down one row,
# carrying its existing formatting along with it).
$ws.Rows.Item(14).Insert()

# New entry: Carlos Salomao / 62f6a9a38f0594d1cedf63bc
$ws.Cells.Item(14,1).Value = "62f6a9a38f0594d1cedf63bc"
$ws.Cells.Item(14,2).Value = "Carlos Salomão"
$ws.Cells.Item(14,3).Value = 370711107.10000002
$ws.Cells.Item(14,4).Value = 13
$ws.Cells.Item(14,5).Value = "l6qv2mk6"

# Renumber the "posicao" column for every row pushed down by the insert
For ($r = 15; $r -le 35; $r++) {
    $ws.Cells.Item($r, 4).Value = $r - 1
}

# Refreshed point totals (column C) for the whole table, as re-synced from source
$ws.Cells.Item(2,3).Value = 2779715780
$ws.Cells.Item(3,3).Value = 2364740766
$ws.Cells.Item(4,3).Value = 2354120787
$ws.Cells.Item(5,3).Value = 1647063840
$ws.Cells.Item(6,3).Value = 1202672017
$ws.Cells.Item(7,3).Value = 1153129978
$ws.Cells.Item(8,3).Value = 968327492.70000005
$ws.Cells.Item(9,3).Value = 648955510.79999995
$ws.Cells.Item(10,3).Value = 590539773.60000002
$ws.Cells.Item(11,3).Value = 579670695.79999995
$ws.Cells.Item(12,3).Value = 450714236.80000001
$ws.Cells.Item(13,3).Value = 411764926.19999999
$ws.Cells.Item(14,3).Value = 370711107.10000002
$ws.Cells.Item(15,3).Value = 289204784.39999998
$ws.Cells.Item(16,3).Value = 252411594.59999999
$ws.Cells.Item(17,3).Value = 233891214.59999999
$ws.Cells.Item(18,3).Value = 205750999.59999999
$ws.Cells.Item(19,3).Value = 138072167.30000001
$ws.Cells.Item(20,3).Value = 114976184.5
$ws.Cells.Item(21,3).Value = 108158853.3
$ws.Cells.Item(22,3).Value = 103098592.90000001
$ws.Cells.Item(23,3).Value = 102258845.8
$ws.Cells.Item(24,3).Value = 66382476.869999997
$ws.Cells.Item(25,3).Value = 56256615.109999999
$ws.Cells.Item(26,3).Value = 53401230.170000002
$ws.Cells.Item(27,3).Value = 52439319.640000001
$ws.Cells.Item(28,3).Value = 51238979.600000001
$ws.Cells.Item(29,3).Value = 45036484.109999999
$ws.Cells.Item(30,3).Value = 33605904.380000003
$ws.Cells.Item(31,3).Value = 31856096.149999999
$ws.Cells.Item(32,3).Value = 29982580.699999999
$ws.Cells.Item(33,3).Value = 20157133.420000002
$ws.Cells.Item(34,3).Value = 7540281.057
$ws.Cells.Item(35,3).Value = 2841504.14

# Row 35 (Jonathan Soares) was appended fresh rather than carried down from the
# shifted block, so it loses the point-column number formatting the others kept.
$ws.Range("A35:E35").ClearFormats()

# Keep the dimension / selection in sync with the extra row
$ws.Range("A2:E35").Select()
